$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44076, 'Hospital Universitario Miguel Servet', 108, 22, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Servcio Aragonés de Salud - 20200901 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44076, 'Hospital Clínico Universitario', 160, 21, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Servcio Aragonés de Salud - 20200901 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44076, 'Hospital Royo Villanova', 40, 5, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Servcio Aragonés de Salud - 20200901 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44076, 'Hospital Nuestra Señora de Gracia', 12, 0, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Servcio Aragonés de Salud - 20200901 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44076, 'Hospital General de la Defensa', 25, 3, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Servcio Aragonés de Salud - 20200901 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44076, 'Hospital Obispo Polanco', 20, 3, 'Teruel', 'Teruel', 44216, 'Fuente Servcio Aragonés de Salud - 20200901 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44076, 'Hospital de Alcañiz', 14, 0, 'Alcañiz', 'Teruel', 44216, 'Fuente Servcio Aragonés de Salud - 20200901 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44076, 'Hospital de Barbastro', 14, 3, 'Barbastro', 'Huesca', 22125, 'Fuente Servcio Aragonés de Salud - 20200901 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44076, 'Hospital San Jorge', 39, 8, 'Huesca', 'Huesca', 22125, 'Fuente Servcio Aragonés de Salud - 20200901 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44076, 'Hospital Sagrado Corazón', 0, 0, 'Huesca', 'Huesca', 22125, 'Fuente Servcio Aragonés de Salud - 20200901 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44076, 'Hospital Ernest Lluch', 5, 0, 'Calatayud', 'Zaragoza', 50297, 'Fuente Servcio Aragonés de Salud - 20200901 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44076, 'Hospital San José', 1, 0, 'Teruel', 'Teruel', 44216, 'Fuente Servcio Aragonés de Salud - 20200901 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44076, 'Hospital Ejea – Cinco Villas', 2, 0, 'Ejea de los Caballeros', 'Zaragoza', 50297, 'Fuente Servcio Aragonés de Salud - 20200901 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44076, 'MAZ', 3, 0, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Servcio Aragonés de Salud - 20200901 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44076, 'Hospital Viamed Montecanal', 1, 0, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Servcio Aragonés de Salud - 20200901 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44076, 'Clínica Montpellier', 6, 0, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Servcio Aragonés de Salud - 20200901 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44076, 'Hospital Quirón', 6, 2, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Servcio Aragonés de Salud - 20200901 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44076, 'Hospital San Juan de Dios de Zaragoza', 30, 0, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Servcio Aragonés de Salud - 20200901 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44076, 'Clínica Viamed Santiago', 2, 0, 'Huesca', 'Huesca', 22125, 'Fuente Servcio Aragonés de Salud - 20200901 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44076, 'Clínica El Pilar', 1, 0, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Servcio Aragonés de Salud - 20200901 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44076, 'C.S. Fraga - Bajo Cinca', 6, 0, 'Fraga', 'Huesca', 22125, 'Fuente Servcio Aragonés de Salud - 20200901 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44076, 'Hospital de Jaca', 0, 0, 'Jaca', 'Huesca', 22125, 'Fuente Servcio Aragonés de Salud - 20200901 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44077, 'Hospital Universitario Miguel Servet', 101, 22, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Servcio Aragonés de Salud - 20200902 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44077, 'Hospital Clínico Universitario', 154, 21, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Servcio Aragonés de Salud - 20200902 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44077, 'Hospital Royo Villanova', 41, 5, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Servcio Aragonés de Salud - 20200902 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44077, 'Hospital Nuestra Señora de Gracia', 11, 0, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Servcio Aragonés de Salud - 20200902 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44077, 'Hospital General de la Defensa', 20, 2, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Servcio Aragonés de Salud - 20200902 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44077, 'Hospital Obispo Polanco', 22, 3, 'Teruel', 'Teruel', 44216, 'Fuente Servcio Aragonés de Salud - 20200902 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44077, 'Hospital de Alcañiz', 16, 0, 'Alcañiz', 'Teruel', 44216, 'Fuente Servcio Aragonés de Salud - 20200902 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44077, 'Hospital de Barbastro', 12, 3, 'Barbastro', 'Huesca', 22125, 'Fuente Servcio Aragonés de Salud - 20200902 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44077, 'Hospital San Jorge', 35, 8, 'Huesca', 'Huesca', 22125, 'Fuente Servcio Aragonés de Salud - 20200902 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44077, 'Hospital Sagrado Corazón', 0, 0, 'Huesca', 'Huesca', 22125, 'Fuente Servcio Aragonés de Salud - 20200902 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44077, 'Hospital Ernest Lluch', 3, 0, 'Calatayud', 'Zaragoza', 50297, 'Fuente Servcio Aragonés de Salud - 20200902 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44077, 'Hospital San José', 1, 0, 'Teruel', 'Teruel', 44216, 'Fuente Servcio Aragonés de Salud - 20200902 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44077, 'Hospital Ejea – Cinco Villas', 7, 0, 'Ejea de los Caballeros', 'Zaragoza', 50297, 'Fuente Servcio Aragonés de Salud - 20200902 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44077, 'MAZ', 3, 0, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Servcio Aragonés de Salud - 20200902 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44077, 'Hospital Viamed Montecanal', 0, 0, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Servcio Aragonés de Salud - 20200902 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44077, 'Clínica Montpellier', 7, 0, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Servcio Aragonés de Salud - 20200902 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44077, 'Hospital Quirón', 5, 2, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Servcio Aragonés de Salud - 20200902 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44077, 'Hospital San Juan de Dios de Zaragoza', 32, 0, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Servcio Aragonés de Salud - 20200902 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44077, 'Clínica Viamed Santiago', 2, 0, 'Huesca', 'Huesca', 22125, 'Fuente Servcio Aragonés de Salud - 20200902 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44077, 'Clínica El Pilar', 1, 0, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Servcio Aragonés de Salud - 20200902 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44077, 'C.S. Fraga - Bajo Cinca', 6, 0, 'Fraga', 'Huesca', 22125, 'Fuente Servcio Aragonés de Salud - 20200902 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
    @(44077, 'Hospital de Jaca', 0, 0, 'Jaca', 'Huesca', 22125, 'Fuente Servcio Aragonés de Salud - 20200902 COVID-19 SITUACION HOSPITALIZACION.xlsx'),
)

$startRow = 2835
$lastExistingRow = 2834

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $srcRange = $ws.Range("A" + $lastExistingRow + ":H" + $lastExistingRow)
    $dstRange = $ws.Range("A" + $row + ":H" + $row)
    $srcRange.Copy($dstRange)

    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
    $ws.Cells.Item($row, 7).Value = $values[6]
    $ws.Cells.Item($row, 8).Value = $values[7]
}

$lastRow = $startRow + $data.Count - 1
$selRange = $ws.Range("A" + ($lastRow - 20) + ":A" + $lastRow)
$selRange.Select() | Out-Null
